$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 46, shifting the existing test rows (old 46-50) down to 47-51.
$ws.Rows("46:46").Insert()

# Populate the new row 46 with a test case for:
# "prevent product to be linked to variation child"
$ws.Range("B46").Value = 42
$ws.Range("C46").Value = 28
$ws.Range("D46").Value = "This tries to link variation to child"
$ws.Range("F46").Value = "Color/Black"
$ws.Range("G46").Value = "Size/XS"
$ws.Range("I46").Value = 12
$ws.Range("J46").Value = 100
$ws.Range("K46").Value = "Test Category"
$ws.Range("L46").Value = "Test Category"
$ws.Range("M46").Value = "shirt1.jpeg"
$ws.Range("N46").Value = "shirt2.jpeg,shirt3.jpeg"

# Selection, as left by the author after editing the sheet.
$ws.Range("B52").Select()

$wb.Save()
